$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 268. This shifts the existing rows
# 268-277 down to 269-278 (values + formats move with them), matching
# the diff's reordering of the trailing block of records.
$ws.Rows(268).Insert()

# Populate the newly inserted row 268 with the new weekly record.
$ws.Range("A268").Value = 8
$ws.Range("B268").Value = "Terminal La Palmera de La Serena"
$ws.Range("C268").Value = "Coquimbo"
$ws.Range("D268").Value = 44939
$ws.Range("E268").Value = 4
$ws.Range("F268").Value = 100112037
$ws.Range("G268").Value = "Cebollín"
$ws.Range("H268").Value = "Sin especificar"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 1160
$ws.Range("K268").Value = 1200
$ws.Range("L268").Value = 1400
$ws.Range("M268").Value = 1300
$ws.Range("N268").Value = "$/paquete 6 unidades"
$ws.Range("O268").Value = "Provincia del Elquí"
$ws.Range("P268").Value = 217
$ws.Range("Q268").Value = 6
$ws.Range("R268").Value = "Hortaliza"
